$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (4th column), shifting
# existing D:G columns (modality0, modality0.source, modality1,
# modality1.source) to E:H.
$ws.Columns("D").Insert()

# Match the bestFit width (13 characters) used by the neighboring columns.
$ws.Columns("D").ColumnWidth = 12.166666666666666

# Header for the new column
$ws.Range("D1").Value = "metadata_dir"

# Fill in metadata_dir values for rows 2 and 3 (row 4 is left blank,
# matching the source data that only provided a metadata_dir for the
# first two data rows)
$ws.Range("D2").Value = "/allen/aind/stage/fake/metadata_dir"
$ws.Range("D3").Value = "/allen/aind/stage/fake/Config"
